$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.328.06'
$ws.Range("E2").Value = '  -0.72%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.711.48'
$ws.Range("E3").Value = '  -0.64%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.27'
$ws.Range("E5").Value = '  -0.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5310'
$ws.Range("E6").Value = '  -1.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06701'
$ws.Range("E8").Value = '  +1.46%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2664'
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.85'
$ws.Range("E10").Value = '  -3.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07657'
$ws.Range("E11").Value = '  -0.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.509'
$ws.Range("E12").Value = '  -2.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.947.55'
$ws.Range("E13").Value = '  -0.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.709.32'
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5826'
$ws.Range("E15").Value = '  -0.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8246'
$ws.Range("E16").Value = '  -0.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.17'
$ws.Range("E17").Value = '  +0.32%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '27.343.95'
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '224.95'
$ws.Range("E19").Value = '  +2.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.629'
$ws.Range("E21").Value = '  -2.08%  '
$ws.Range("E22").Value = '  -2.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.011'
$ws.Range("E23").Value = '  -1.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.003'
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.49'
$ws.Range("E25").Value = '  -2.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.692'
$ws.Range("E26").Value = '  -2.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1207'
$ws.Range("E27").Value = '  -2.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.239'
$ws.Range("E28").Value = '  -2.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '16.27'
$ws.Range("E29").Value = '  -2.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05366'
$ws.Range("E30").Value = '  -3.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.292'
$ws.Range("E31").Value = '  -0.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.487'
$ws.Range("E32").Value = '  -2.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.439'
$ws.Range("E33").Value = '  -0.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.637'
$ws.Range("E34").Value = '  -1.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.866'
$ws.Range("E35").Value = '  +1.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9505'
$ws.Range("E36").Value = '  -1.15%  '
$ws.Range("E37").Value = '  -1.50%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5850'
$ws.Range("E38").Value = '  -2.01%  '
$ws.Range("E39").Value = '  -0.72%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.085.71'
$ws.Range("E40").Value = '  +3.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.793'
$ws.Range("E41").Value = '  -2.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.003'
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8390'
$ws.Range("E43").Value = '  -1.88%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.72'
$ws.Range("E44").Value = '  -0.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.854.78'
$ws.Range("E45").Value = '  -0.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈115'
$ws.Range("E46").Value = '  -0.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '57.72'
$ws.Range("E47").Value = '  -2.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4526'
$ws.Range("E48").Value = '  +2.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.004'
$ws.Range("E49").Value = '  +0.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.112'
$ws.Range("E50").Value = '  -1.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05227'
$ws.Range("E51").Value = '  -0.33%  '
